$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38e095848ec25b80ff89be18825606c120ad93542e03c2755eb84ee90edfe5b5"
$ws.Range("D3").Value = "245d2c837d0ca60011818cdb9aded89abfb4d4f5eb94148440124e43d291921d"
$ws.Range("D4").Value = "1a69917e78ae842865dfb5571a62d65494a59aac74a48f31f3bb799425cdf535"
$ws.Range("D5").Value = "915b10f800bf6b1db00913183b44f454e0b4daa910915de907369297ff61dcda"
$ws.Range("D6").Value = "6b11204cfb0ef5ded49c9bad3ecf995a34c45a0c7a33ea388da12f409b750f73"
